# Apply the "Fixed angle profile connector" change:
# Add a new row of data (Chassis_Bottom_Height = 80 mm) into row 9
# of the Chassis-Dimensions sheet, and move the active selection to B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Chassis_Bottom_Height"
$ws.Range("B9").Value = 80
$ws.Range("C9").Value = "mm"

# Match the numeric formatting used by the other "Values" cells (e.g. B8)
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B9").Value = 80

$ws.Range("B9").Select()
